# "Added reports folder and a new dump"
# This workbook is a weekly "Outstanding Results" report. The edit:
#   - Refreshes the generated/sent timestamp in K6.
#   - Inserts a new "Date Sent" header column (shifting "Death Register No."
#     into the column that used to read "Number of Samples Sent").
#   - Clears out the two sample detail rows (12 & 13) that held the old
#     dump's example/leftover figures, leaving the rows blank for the new
#     dump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp in the merged K6:N6 cell (when this report snapshot was produced)
$ws.Range("K6").Value = 41439.3249241088

# Header row 11 - add "Date Sent" as the first column and shift the
# "Death Register No." label into the next column
$ws.Range("A11").Value = "Date Sent"
$ws.Range("C11").Value = "Death Register No."

# Row 12 - clear the stale sample data
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("I12").Value = ""

# Row 13 - clear the stale sample data
$ws.Range("A13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("I13").Value = ""
